# edit.ps1 - applies the "More comments for Haiquan and Qike" revision
#
# Summary of content changes reproduced here:
#  1. Comment ("Is it possible to use more than two with this method? ...")
#     originally stored as two runs split by a stray "_GoBack" bookmark;
#     normalize it to a single run with the bookmark removed.
#  2. Comment ("I like the figure but the plot points and text should be
#     larger. ...") gets a new trailing sentence appended:
#     "  Also, there are no axis labels"
#  3. Fix the typo "caner" -> "cancer" in the "Problems and alternatives"
#     paragraph via a tracked insertion of the letter "c" (author
#     "Dominic LaRoche [2]"), leaving a "_GoBack" bookmark right after the
#     newly inserted letter (Word drops its last-edit marker there).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the split comment run and drop the embedded _GoBack bookmark.
#    (Comments.Item index is ordered by the position of their anchor in
#    the document; this is the "Is it possible to use more than two..."
#    comment anchored on the word "two".)
# ---------------------------------------------------------------------
$fixedComment = $null
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $candidate = $d.Comments.Item($i)
    if ($candidate.Range.Text -like "Is it possible to use more than two*") {
        $fixedComment = $candidate
    }
}

if ($fixedComment -ne $null) {
    $mergedText = "Is it possible to use more than two with this method?  That would be something worth mentioning even if you do not use more than two in your demonstrational analysis."
    # Force a real rewrite even though the visible text ends up identical
    # (a no-op assignment is detected and skipped), which is what actually
    # merges the runs and removes the bookmark sitting between them.
    $fixedComment.Range.Text = "placeholder"
    $fixedComment2 = $d.Comments.Item($fixedComment.Index)
    $fixedComment2.Range.Text = $mergedText
}

# ---------------------------------------------------------------------
# 2) Append a new sentence onto the end of the figure-quality comment.
# ---------------------------------------------------------------------
$figureComment = $null
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $candidate = $d.Comments.Item($i)
    if ($candidate.Range.Text -like "I like the figure*") {
        $figureComment = $candidate
    }
}

if ($figureComment -ne $null) {
    $appendedText = $figureComment.Range.Text + "  Also, there are no axis labels"
    $figureComment.Range.Text = $appendedText
}

# ---------------------------------------------------------------------
# 3) Fix "caner" -> "cancer" with a tracked insertion, then drop the
#    _GoBack bookmark right after the inserted letter.
# ---------------------------------------------------------------------
$word.UserName = "Dominic LaRoche [2]"
$previousTrackRevisions = $d.TrackRevisions
$d.TrackRevisions = $true

$searchRange = $d.Content
$found = $searchRange.Find.Execute("the same caner type.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $splitOffset = $searchRange.Start + "the same can".Length
    $insertionPoint = $d.Range($splitOffset, $splitOffset)
    $insertionPoint.InsertAfter("c")

    $goBackPoint = $d.Range($splitOffset + 1, $splitOffset + 1)
    $d.Bookmarks.Add("_GoBack", $goBackPoint)
}

$d.TrackRevisions = $previousTrackRevisions
